$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new columns before column D to hold the two newest quarters
$ws.Range("D:E").Insert()

# Copy cell formatting (number format / font / style) from column F (the former column D) into new D:E columns
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate new column D and E values for each row
$ws.Range("D7").Value = 43463
$ws.Range("E7").Value = 43372
$ws.Range("D8").Value = 114600
$ws.Range("E8").Value = 290300
$ws.Range("D9").Value = 78500
$ws.Range("E9").Value = 192800
$ws.Range("D10").Value = 36100
$ws.Range("E10").Value = 97500
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 23600
$ws.Range("E14").Value = 700
$ws.Range("D15").Value = 900
$ws.Range("E15").Value = 900
$ws.Range("D17").Value = 153700
$ws.Range("E17").Value = 253100
$ws.Range("D18").Value = -39100
$ws.Range("E18").Value = 37300
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = -34700
$ws.Range("E21").Value = 41400
$ws.Range("D22").Value = 4200
$ws.Range("E22").Value = 4200
$ws.Range("D23").Value = -43200
$ws.Range("E23").Value = 33100
$ws.Range("D24").Value = -4600
$ws.Range("E24").Value = 14500
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -38600
$ws.Range("E26").Value = 18600
$ws.Range("D27").Value = -38600
$ws.Range("E27").Value = 18600
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = -38600
$ws.Range("E33").Value = 18600
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -38600
$ws.Range("E35").Value = 18600
$ws.Range("D38").Value = 43463
$ws.Range("E38").Value = 43372
$ws.Range("D41").Value = 1000
$ws.Range("E41").Value = 7900
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 81900
$ws.Range("E43").Value = 178800
$ws.Range("D44").Value = 90100
$ws.Range("E44").Value = 96000
$ws.Range("D45").Value = 12800
$ws.Range("E45").Value = 14000
$ws.Range("D46").Value = 185800
$ws.Range("E46").Value = 296800
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 31900
$ws.Range("E48").Value = 31700
$ws.Range("D49").Value = 52700
$ws.Range("E49").Value = 76500
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 300
$ws.Range("E52").Value = 2000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 270700
$ws.Range("E54").Value = 407000
$ws.Range("D57").Value = 41300
$ws.Range("E57").Value = 43200
$ws.Range("D58").Value = 30400
$ws.Range("E58").Value = 90500
$ws.Range("D59").Value = 27000
$ws.Range("E59").Value = 38400
$ws.Range("D60").Value = 98600
$ws.Range("E60").Value = 172100
$ws.Range("D61").Value = 103600
$ws.Range("E61").Value = 128800
$ws.Range("D62").Value = 1100
$ws.Range("E62").Value = 600
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 203300
$ws.Range("E66").Value = 301500
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -55600
$ws.Range("E72").Value = -17000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 67400
$ws.Range("E76").Value = 105500
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43463
$ws.Range("E80").Value = 43372
$ws.Range("D81").Value = -38600
$ws.Range("E81").Value = 18600
$ws.Range("D83").Value = 4300
$ws.Range("E83").Value = 4200
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 83500
$ws.Range("E89").Value = -20100
$ws.Range("D91").Value = -3500
$ws.Range("E91").Value = -2900
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -4400
$ws.Range("E94").Value = -3900
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -86100
$ws.Range("E100").Value = 23300
$ws.Range("D101").Value = 200
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -6900
$ws.Range("E102").Value = -700

# Fix values in shifted columns that differ from a pure shift (explicit corrections from source data)
$ws.Range("F15").Value = 1000
$ws.Range("G15").Value = 1000
$ws.Range("F91").Value = -2600
$ws.Range("G91").Value = -3400
$ws.Range("H91").Value = -3100
$ws.Range("I91").Value = -3500
$ws.Range("J91").Value = -4900

Write-Host "Edit complete"
